$wb = $excel.ActiveWorkbook

# --- objective sheet ---
$wsObjective = $wb.Worksheets.Item("objective")
$wsObjective.Range("A2").Value = 55.542173157718125

# --- mass_balance sheet ---
$wsMassBalance = $wb.Worksheets.Item("mass_balance")
$wsMassBalance.Range("D2").Value = -0.0
$wsMassBalance.Range("D3").Value = -21.0
$wsMassBalance.Range("D4").Value = -0.0
$wsMassBalance.Range("D5").Value = -0.3855432894295314
$wsMassBalance.Range("D6").Value = -11.0
$wsMassBalance.Range("D7").Value = -0.0
$wsMassBalance.Range("D8").Value = -0.3855432894295314
$wsMassBalance.Range("D9").Value = -0.0
$wsMassBalance.Range("D10").Value = -0.0
$wsMassBalance.Range("D11").Value = -0.0
$wsMassBalance.Range("D12").Value = -0.0
$wsMassBalance.Range("D13").Value = -0.0
